# Apply "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) inside specific bullet
# points, matching the target OOXML diff. Numbers/measurements get a
# dedicated run with Bold + the brand color (#2C3E50) while the
# surrounding text stays in plain, unformatted runs.

$d = $word.ActiveDocument

# Color used for highlighted metrics: RGB(0x2C, 0x3E, 0x50) encoded as
# Word's 0xBBGGRR integer (the format Font.Color / WdColor expects).
$metricColor = 0x2C + (0x3E * 256) + (0x50 * 65536)   # 5258796

function Find-ParagraphContaining($anchorText, $excludeText) {
    # Returns the Word Paragraph object whose text contains $anchorText
    # (first match, scanning top to bottom) and, if $excludeText is
    # non-empty, does NOT contain $excludeText. Using a content anchor
    # instead of a hardcoded paragraph index keeps this resilient to
    # minor structural drift elsewhere in the document.
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Contains($anchorText)) {
            if ([string]::IsNullOrEmpty($excludeText) -or -not $t.Contains($excludeText)) {
                return $p
            }
        }
    }
    throw "Paragraph not found containing: $anchorText"
}

function Set-MetricFormat($range, $searchText) {
    # Finds the first occurrence of $searchText inside $range's current
    # text and applies bold + the metric color to just that substring,
    # leaving the rest of the paragraph's formatting untouched.
    $full = $range.Text
    $idx = $full.IndexOf($searchText)
    if ($idx -lt 0) {
        throw "Segment not found: $searchText"
    }
    $segStart = $range.Start + $idx
    $segEnd = $segStart + $searchText.Length
    $seg = $d.Range($segStart, $segEnd)
    $seg.Font.Bold = 1
    $seg.Font.Color = $metricColor
}

function Format-ParagraphMetrics($anchorText, $excludeText, $segments) {
    $para = Find-ParagraphContaining $anchorText $excludeText
    $range = $para.Range
    foreach ($seg in $segments) {
        Set-MetricFormat $range $seg
    }
}

# 1) "• Discovered systematic race coding errors ... from 23% to 64%"
#    (the long achievement bullet under Siege Analytics, not the short
#    "KEY ACHIEVEMENTS" one, so anchor on the longer, unique phrase)
Format-ParagraphMetrics "developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" "" @("23%", "64%")

# 2) "• Achieved 87% prediction accuracy ... from ±4.2% to ±2.1%"
#    (the long version with the polling-error-margin clause; the ± sign
#    makes this anchor unique to this paragraph)
Format-ParagraphMetrics ([char]0x00B1 + "4.2%") "" @("87%", "71%", ([char]0x00B1 + "4.2%"), ([char]0x00B1 + "2.1%"))

# 3) "• Wrote RFP and analyzed bids from 1,200 vendors ..."
Format-ParagraphMetrics "Wrote RFP and analyzed bids from 1,200 vendors" "" @("1,200")

# 4) "• Created comprehensive meta-analysis framework ... $400M ... $1B+"
Format-ParagraphMetrics "Polling Consortium Database" "" @(("$" + "400M"), ("$" + "1B"))

# 5) "• Algorithm reduced mapping costs by 73.5%, saving ... $4.7M"
Format-ParagraphMetrics "Algorithm reduced mapping costs" "" @("73.5%", ("$" + "4.7M"))

# 6) "• Achieved 87% prediction accuracy for voter turnout vs. industry
#    standard of 71%" (the short version, under KEY ACHIEVEMENTS AND
#    IMPACT bullet list - exclude the long variant's trailing clause to
#    disambiguate from paragraph #2 above)
Format-ParagraphMetrics "industry standard of 71%" "reducing polling error margins" @("87%", "71%")

Write-Output "metrics highlighted"
